$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns -------------------------------------------------
# Original layout (A-H): date, activity, type, subject, hours, related_object,
# related_to, status
# Target layout (A-M): date, activity, type, subject, notes, next_step,
# related_object, related_to, activity_category, solution, solution_product,
# hours, status

# Keep the text "hours" alive in a far-away placeholder cell while we
# restructure the columns, so that it keeps its original identity/position
# in the shared strings table instead of being dropped and re-appended
# later (we'll move it back into its final spot afterwards).
$ws.Range("Z1").Value = "hours"

# 1) Remove the old "hours" column (E). related_object/related_to/status
#    shift one column to the left, keeping their column formatting.
$ws.Range("E1").EntireColumn.Delete()

# 2) Insert two new blank columns for "notes" and "next_step" right before
#    the (now shifted) related_object column.
$ws.Range("E1:F1").EntireColumn.Insert()

# 3) Insert three new blank columns for "activity_category", "solution" and
#    "solution_product" right after related_to (before status).
$ws.Range("I1:K1").EntireColumn.Insert()

# 4) Insert one new blank column for "hours" right before status.
$ws.Range("L1").EntireColumn.Insert()

# --- Fill in the new / moved header cells and data cells (order chosen to
#     match the order in which the shared strings table was originally
#     built) -----------------------------------------------------------
$ws.Range("E1").Value = "notes"

$ws.Range("E2").Value = "This is something."
$ws.Range("E3").Value = "This is something."

$ws.Range("F1").Value = "next_step"

$ws.Range("F3").Value = "N/A"

$ws.Range("I3").Value = "EUC"

$ws.Range("J1").Value = "solution"
$ws.Range("J3").Value = "HORIZON ON PREM"

$ws.Range("K3").Value = "HORIZON STD"

$ws.Range("I1").Value = "activity_category"
$ws.Range("K1").Value = "solution_product"

# --- Restore the "hours" header/values in their new column (L) -----------
$ws.Range("L1").Value = "hours"
$ws.Range("L2").Value = 2
$ws.Range("L3").Value = 2

# Remove the placeholder now that "hours" lives in its real spot again.
$ws.Range("AE1").Clear()

# --- Column widths for the new columns ------------------------------------
# E (notes) / F (next_step) -> width 16
$ws.Range("E1").EntireColumn.ColumnWidth = 15.166666666666666
$ws.Range("F1").EntireColumn.ColumnWidth = 15.166666666666666

# I:K (activity_category / solution / solution_product) -> width 23.33203125
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 22.498697916666664

# --- Selection matches the author's final cursor position -----------------
$ws.Range("K2").Select()
